$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Update the Status column text (shared by Overview + per-locale sheets)
# ---------------------------------------------------------------------------
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# helper: rebuild a locale sheet's hyperlinks, inserting the new
# "Latest Target File" (F) / "Latest Handback File" (G) columns, and
# refresh the "Latest Handback DateTime" (H) column.
# ---------------------------------------------------------------------------
function Update-LocaleSheet($ws, $handbackTime) {
    # capture the existing hyperlink addresses before wiping them out
    $addrA2 = ""
    $addrD2 = ""
    $addrA3 = ""
    $addrD3 = ""
    foreach ($hl in $ws.Hyperlinks) {
        $rngAddr = $hl.Range.Address()
        if ($rngAddr -eq '$A$2') { $addrA2 = $hl.Address }
        if ($rngAddr -eq '$D$2') { $addrD2 = $hl.Address }
        if ($rngAddr -eq '$A$3') { $addrA3 = $hl.Address }
        if ($rngAddr -eq '$D$3') { $addrD3 = $hl.Address }
    }

    $mdName = $ws.Range("A2").Value2
    $xlfName = $ws.Range("D2").Value2
    $mdName3 = $ws.Range("A3").Value2

    # drop every hyperlink on the sheet so we can recreate them in the
    # desired left-to-right / top-to-bottom order
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $addrA2, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("D2"), $addrD2, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("F2"), $addrA2, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("G2"), $addrD2, "", "", $xlfName)

    $ws.Hyperlinks.Add($ws.Range("A3"), $addrA3, "", "", $mdName3)
    $ws.Hyperlinks.Add($ws.Range("D3"), $addrD3, "", "", $xlfName)
    $ws.Hyperlinks.Add($ws.Range("F3"), $addrA2, "", "", $mdName)
    $ws.Hyperlinks.Add($ws.Range("G3"), $addrD2, "", "", $xlfName)

    # give the two new columns the same "HyperLink" look as the rest of
    # the hyperlinked cells (underlined, cornflower blue)
    foreach ($addr in @("F2", "G2", "F3", "G3")) {
        $c = $ws.Range($addr)
        $c.Font.Name = "Calibri"
        $c.Font.Size = 11
        $c.Font.Color = 15570276
        $c.Font.Underline = 2
    }

    # refresh the "Latest Handback DateTime" column
    $ws.Range("H2").Value = $handbackTime
    $ws.Range("H3").Value = $handbackTime
}

Update-LocaleSheet $wsZh "2016-03-20 17:26:13"
Update-LocaleSheet $wsDe "2016-03-20 17:26:26"
